$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new product rows at the end of the sheet (rows 126-127),
# extending the used range from A1:C125 to A1:C127.

$ws.Cells.Item(126, 1).Value = 92007776
$ws.Cells.Item(126, 2).Value = "Jogo de Brocas para Metal e Madeira 19 Peças Dexter"
$ws.Cells.Item(126, 3).Value = 59.9

$ws.Cells.Item(127, 1).Value = 89235783
$ws.Cells.Item(127, 2).Value = "Luva Nylon Látex Maxigrip Pro M Danny"
$ws.Cells.Item(127, 3).Value = 17.49
